$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("STM32")

# Add the four new URF sensors (5-8): echo pin label in column A,
# trig pin label in column F, for their respective rows.
# New shared strings must be created in "echo" order first, then "trig"
# order, to match the canonical string table layout.
$ws.Range("A16").Value = "URF5 echo"
$ws.Range("A17").Value = "URF6 echo"
$ws.Range("A19").Value = "URF7 echo"
$ws.Range("A20").Value = "URF8 echo"

$ws.Range("F16").Value = "URF5 trig"
$ws.Range("F17").Value = "URF6 trig"
$ws.Range("F19").Value = "URF7 trig"
$ws.Range("F20").Value = "URF8 trig"

# Update the view so the scrolled-to cell / selection matches where the
# edits were made (topLeftCell A10 -> A7, selection E17 -> F16).
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("F16").Select()
